# Updating with Caroline Krafft workshop
#
# On the "2024 - Fall" sheet, row 9 is the Caroline Krafft "Surveys Design"
# workshop (Oct. 23rd). It previously had no recorded attendance (Faculty/
# Alumni = 0, Grad/Undergrad columns empty). Record the attendance: 4 in
# column G (Faculty, Alumni) and 9 in column H (Grad, Undergrad) -> I9
# (total) becomes 13, which cascades through the running-total column K and
# the summary rows below, as well as into the "Attendance Descriptives"
# sheet that pulls from this sheet's K column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")

$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 9

# Reflect the editor's final cursor position (cell H10) that was recorded
# in the saved file after entering the H9 value.
$ws.Activate()
$ws.Range("H10").Select()
